$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''295.94'
$ws.Range('E2').Value = '''3.34%'
$ws.Range('G2').Value = '''17'
$ws.Range('D3').Value = '''41.31'
$ws.Range('E3').Value = '''2.93%'
$ws.Range('G3').Value = '''17'
$ws.Range('D4').Value = '''5.024'
$ws.Range('E4').Value = '''-0.05%'
$ws.Range('G4').Value = '''17'
$ws.Range('D5').Value = '''0.07457'
$ws.Range('E5').Value = '''2.32%'
$ws.Range('G5').Value = '''17'
$ws.Range('D6').Value = '''1.576'
$ws.Range('E6').Value = '''3.28%'
$ws.Range('G6').Value = '''17'
$ws.Range('D7').Value = '''0.9269'
$ws.Range('E7').Value = '''1.23%'
$ws.Range('G7').Value = '''17'
$ws.Range('G8').Value = '''17'
$ws.Range('D9').Value = '''0.1184'
$ws.Range('E9').Value = '''-1.04%'
$ws.Range('G9').Value = '''17'
$ws.Range('D10').Value = '''0.1806'
$ws.Range('E10').Value = '''5.27%'
$ws.Range('G10').Value = '''17'
$ws.Range('D11').Value = '''0.08820'
$ws.Range('E11').Value = '''1.96%'
$ws.Range('G11').Value = '''17'
$ws.Range('D12').Value = '''0.04188'
$ws.Range('E12').Value = '''-0.15%'
$ws.Range('G12').Value = '''17'
$ws.Range('D13').Value = '''0.1049'
$ws.Range('E13').Value = '''-0.25%'
$ws.Range('G13').Value = '''17'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '''0.001273'
$ws.Range('E14').Value = '''0.38%'
$ws.Range('G14').Value = '''17'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = '''0.005889'
$ws.Range('E15').Value = '''-1.69%'
$ws.Range('G15').Value = '''17'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = '''3.355'
$ws.Range('E16').Value = '''-1.33%'
$ws.Range('G16').Value = '''17'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = '''4.353'
$ws.Range('E17').Value = '''0.98%'
$ws.Range('G17').Value = '''17'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D18').Value = '''0.3289'
$ws.Range('E18').Value = '''0.15%'
$ws.Range('G18').Value = '''17'
$ws.Range('B19').Value = 'MCDex'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D19').Value = '''7.862'
$ws.Range('E19').Value = '''-0.36%'
$ws.Range('G19').Value = '''17'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').Value = '''0.1412'
$ws.Range('E20').Value = '''5.11%'
$ws.Range('G20').Value = '''17'
$ws.Range('B21').Value = 'ZBToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D21').Value = '''0.2972'
$ws.Range('E21').Value = '''2.92%'
$ws.Range('G21').Value = '''17'
$ws.Range('D22').Value = '''0.04025'
$ws.Range('E22').Value = '''4.83%'
$ws.Range('G22').Value = '''17'
$ws.Range('D23').Value = '''0.001265'
$ws.Range('E23').Value = '''-0.40%'
$ws.Range('G23').Value = '''17'
$ws.Range('D24').Value = '''0.003862'
$ws.Range('E24').Value = '''2.64%'
$ws.Range('G24').Value = '''17'
$ws.Range('D25').Value = '''0.0001230'
$ws.Range('E25').Value = '''-4.15%'
$ws.Range('G25').Value = '''17'
$ws.Range('D26').Value = '''0.0003720'
$ws.Range('E26').Value = '''-0.30%'
$ws.Range('G26').Value = '''17'
$ws.Range('G27').Value = '''17'
$ws.Range('G28').Value = '''17'
$ws.Range('G29').Value = '''17'
$ws.Range('G30').Value = '''17'
$ws.Range('G31').Value = '''17'
$ws.Range('G32').Value = '''17'
$ws.Range('G33').Value = '''17'
$ws.Range('G34').Value = '''17'
$ws.Range('G35').Value = '''17'
$ws.Range('G36').Value = '''17'
$ws.Range('G37').Value = '''17'
$ws.Range('D38').Value = '''0.02390'
$ws.Range('E38').Value = '''3.81%'
$ws.Range('G38').Value = '''17'
$ws.Range('D39').Value = '''0.05200'
$ws.Range('E39').Value = '''4.61%'
$ws.Range('G39').Value = '''17'
$ws.Range('D40').Value = '''0.006236'
$ws.Range('E40').Value = '''-7.92%'
$ws.Range('G40').Value = '''17'
$ws.Range('D41').Value = '''0.007786'
$ws.Range('E41').Value = '''1.13%'
$ws.Range('G41').Value = '''17'
$ws.Range('D42').Value = '''0.1313'
$ws.Range('E42').Value = '''3.82%'
$ws.Range('G42').Value = '''17'
$ws.Range('D43').Value = '''0.007374'
$ws.Range('E43').Value = '''-0.21%'
$ws.Range('G43').Value = '''17'
$ws.Range('D44').Value = '''0.007799'
$ws.Range('E44').Value = '''4.73%'
$ws.Range('G44').Value = '''17'
$ws.Range('D45').Value = '''0.3215'
$ws.Range('E45').Value = '''4.37%'
$ws.Range('G45').Value = '''17'
$ws.Range('D46').Value = '''0.00006243'
$ws.Range('E46').Value = '''-2.38%'
$ws.Range('G46').Value = '''17'
$ws.Range('D47').Value = '''0.00000000750'
$ws.Range('E47').Value = '''-0.30%'
$ws.Range('G47').Value = '''17'
$ws.Range('D48').Value = '''0.04605'
$ws.Range('E48').Value = '''-81.71%'
$ws.Range('G48').Value = '''17'
$ws.Range('D49').Value = '''0.004198'
$ws.Range('E49').Value = '''-0.08%'
$ws.Range('G49').Value = '''17'
$ws.Range('D50').Value = '''0.00002099'
$ws.Range('E50').Value = '''-0.30%'
$ws.Range('G50').Value = '''17'
$ws.Range('D51').Value = '''0.0001999'
$ws.Range('E51').Value = '''-0.30%'
$ws.Range('G51').Value = '''17'
